$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2026-01-10 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-11 Sunday", 2)

# Update the 25 multiplication problems laid out in the single 20-row x 5-col
# table (content only lives in rows 1, 5, 10, 15, 20 - the rest are spacer
# rows). Addressing cells directly by (row, col) avoids any ambiguity from
# duplicate text values appearing before/after the edit.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "107×5=535"
$t.Cell(1, 2).Range.Text  = "927×8=7416"
$t.Cell(1, 3).Range.Text  = "427×3=1281"
$t.Cell(1, 4).Range.Text  = "942×2=1884"
$t.Cell(1, 5).Range.Text  = "983×4=3932"

$t.Cell(5, 1).Range.Text  = "930×8=7440"
$t.Cell(5, 2).Range.Text  = "511×5=2555"
$t.Cell(5, 3).Range.Text  = "228×6=1368"
$t.Cell(5, 4).Range.Text  = "416×6=2496"
$t.Cell(5, 5).Range.Text  = "129×4=516"

$t.Cell(10, 1).Range.Text = "380×3=1140"
$t.Cell(10, 2).Range.Text = "230×4=920"
$t.Cell(10, 3).Range.Text = "347×4=1388"
$t.Cell(10, 4).Range.Text = "804×3=2412"
$t.Cell(10, 5).Range.Text = "616×4=2464"

$t.Cell(15, 1).Range.Text = "282×2=564"
$t.Cell(15, 2).Range.Text = "290×3=870"
$t.Cell(15, 3).Range.Text = "378×8=3024"
$t.Cell(15, 4).Range.Text = "780×6=4680"
$t.Cell(15, 5).Range.Text = "846×5=4230"

$t.Cell(20, 1).Range.Text = "405×8=3240"
$t.Cell(20, 2).Range.Text = "542×9=4878"
$t.Cell(20, 3).Range.Text = "432×9=3888"
$t.Cell(20, 4).Range.Text = "704×4=2816"
$t.Cell(20, 5).Range.Text = "276×3=828"
